$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 0.007
$ws.Range("C5").Value = 0.027
$ws.Range("C10").Value = 0.601
$ws.Range("C11").Value = 0.949
$ws.Range("C12").Value = 0.145
$ws.Range("C13").Value = 0.486

$wb.Save()
